$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (shared strings will be regenerated accordingly)
$ws.Range("A2").Value = "P2210-00003"
$ws.Range("C2").Value = "Phiếu nhập tháng 10"
$ws.Range("D2").Value = "11-10-2022 00:00:00"
$ws.Range("E2").Value = "1.100.000 VND"
$ws.Range("F2").Value = "Phiếu nhập tháng 10"

# Adjust column widths to match new content (stored width = ColumnWidth + 5/7)
$ws.Columns.Item(3).ColumnWidth = 22.2857142857143
$ws.Columns.Item(6).ColumnWidth = 22.2857142857143
